# Add a new "Transfer" user story with its tasks, and a new "provision"
# (fee calculation) user story with its tasks, to the estimation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 29: new user story header ("WalletTransfer")
$ws.Range("A29").Value = "Kao korisnik potrebno je da mogu da vrsim transfer novca sa svog na neki drugi novcanik u sistemu"
$ws.Range("A29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 30

# Rows 30-33: tasks for the Transfer user story
$ws.Range("B30").Value = "Dodavanje Transfer metode na WalletService"
$ws.Range("C30").Value = 20

$ws.Range("B31").Value = "Implementacija testova za Transfer"
$ws.Range("C31").Value = 30

$ws.Range("B32").Value = "Dodavanje rute za Transfer sredstava u WalletController"
$ws.Range("C32").Value = 10

$ws.Range("B33").Value = "Dodavanje stranice za Transfer sredstava u MVC aplikaciju"
$ws.Range("C33").Value = 10

# Row 34: new user story header (provision / fee calculation), plain style
$ws.Range("A34").Value = "Sistem treba da racuna proviziju u slucaju transfera novca izmedju dva novcanika"

# Rows 35-37: tasks for the provision user story
$ws.Range("B35").Value = "Implementacija servisa za racunanje provizije"
$ws.Range("C35").Value = 60

$ws.Range("B36").Value = "Implementacija testova za racunanje provizije"
$ws.Range("C36").Value = 30

$ws.Range("B37").Value = "Prikaz provizije na transfer stranici"
$ws.Range("C37").Value = 10

# Restore the selection to match the author's final cursor position
$ws.Range("B36").Select()
